$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new log entry in row 82: date value and activity text (new shared string)
$ws.Cells.Item(82, 1).Value = 41407
$ws.Cells.Item(82, 2).Value = "Added benchmark of fastest matrix alg"

# Update the active selection to reflect the new entry location
$ws.Range("B88").Select()
